$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.021.22"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.32"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.77"
$ws.Range("E5").Value = "  -2.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4587"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.75"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07970"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9916"
$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("E12").Value = "  -3.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.899.96"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.906"
$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.073"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.34"
$ws.Range("E17").Value = "  -3.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001030"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.042.09"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  -2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.38"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.207"
$ws.Range("E25").Value = "  -2.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.95"
$ws.Range("E26").Value = "  -1.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.53"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.58"
$ws.Range("E28").Value = "  -2.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.098"

$ws.Range("E30").Value = "  -1.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.81"
$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.005"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09329"
$ws.Range("E33").Value = "  -2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.601"
$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.276"
$ws.Range("E36").Value = "  -2.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06054"
$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("E38").Value = "  -3.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.265"
$ws.Range("E39").Value = "  -4.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.173"
$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5770"
$ws.Range("E42").Value = "  -4.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1824"
$ws.Range("E43").Value = "  -4.27%  "

$ws.Range("E44").Value = "  -4.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.261"
$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07523"
$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.262"
$ws.Range("E47").Value = "  +4.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.05"
$ws.Range("E48").Value = "  -2.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5452"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.896"
$ws.Range("E50").Value = "  -4.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.13"
$ws.Range("E51").Value = "  -1.75%  "
